$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1424977709956892
$ws.Range("D2").Value = 0.005807487234779529
$ws.Range("E2").Value = 0.1057584759616947
$ws.Range("F2").Value = 0.3718757801769073
$ws.Range("G2").Value = 0.2574595956363908
$ws.Range("H2").Value = 0.3318063431118219
$ws.Range("M2").Value = 0.5878127933305421
$ws.Range("N2").Value = 1.469340516368959
$ws.Range("O2").Value = 1.094333041502466

$ws.Range("B3").Value = 0.1329151663101698
$ws.Range("D3").Value = 0.005159819009612931
$ws.Range("E3").Value = 0.1086794443102992
$ws.Range("F3").Value = 0.3503061484837247
$ws.Range("G3").Value = 0.2357554001262656
$ws.Range("H3").Value = 0.325221279927348
$ws.Range("M3").Value = 0.5138342799021132
$ws.Range("N3").Value = 1.412118234844286
$ws.Range("O3").Value = 1.034226808312042

$ws.Range("B4").Value = 0.1270974125729936
$ws.Range("D4").Value = 0.004762600002866435
$ws.Range("E4").Value = 0.110665057293458
$ws.Range("F4").Value = 0.3373112911113338
$ws.Range("G4").Value = 0.2225899383044663
$ws.Range("H4").Value = 0.3213857305289167
$ws.Range("M4").Value = 0.4682921863596761
$ws.Range("N4").Value = 1.377486344751787
$ws.Range("O4").Value = 0.9980857925205555

$ws.Range("B5").Value = 0.1247434455953282
$ws.Range("D5").Value = 0.004600856605890158
$ws.Range("E5").Value = 0.1115221410048459
$ws.Range("F5").Value = 0.3320782194013532
$ws.Range("G5").Value = 0.2172651167785915
$ws.Range("H5").Value = 0.3198749083758514
$ws.Range("M5").Value = 0.4497044331874633
$ws.Range("N5").Value = 1.363502128123514
$ws.Range("O5").Value = 0.9835496024267059

$ws.Range("B6").Value = 0.1243535931561581
$ws.Range("D6").Value = 0.004574007262860391
$ws.Range("E6").Value = 0.1116673428104189
$ws.Range("F6").Value = 0.3312130379675935
$ws.Range("G6").Value = 0.216383354860028
$ws.Range("H6").Value = 0.3196271890059137
$ws.Range("M6").Value = 0.4466162284019362
$ws.Range("N6").Value = 1.361187886994799
$ws.Range("O6").Value = 0.9811474300563532

$ws.Range("B7").Value = 0.1270655977794064
$ws.Range("D7").Value = 0.004760418144439171
$ws.Range("E7").Value = 0.1106764226615748
$ws.Range("F7").Value = 0.3372404634936643
$ws.Range("G7").Value = 0.2225179634499597
$ws.Range("H7").Value = 0.3213651437619092
$ws.Range("M7").Value = 0.4680416215080356
$ws.Range("N7").Value = 1.377297225357268
$ws.Range("O7").Value = 0.9978889776081701

$ws.Range("B8").Value = 0.1391801169005475
$ws.Range("D8").Value = 0.005584085006280048
$ws.Range("E8").Value = 0.1067255206982551
$ws.Range("F8").Value = 0.3643867845611481
$ws.Range("G8").Value = 0.2499424014263525
$ws.Range("H8").Value = 0.3294926804022396
$ws.Range("M8").Value = 0.5623302377808272
$ws.Range("N8").Value = 1.449507158988155
$ws.Range("O8").Value = 1.073449331504975

$ws.Range("B9").Value = 0.1634521087393352
$ws.Range("D9").Value = 0.007202387369975582
$ws.Range("E9").Value = 0.1005193123292933
$ws.Range("F9").Value = 0.4196070520171418
$ws.Range("G9").Value = 0.3050143363936542
$ws.Range("H9").Value = 0.3470816658063711
$ws.Range("M9").Value = 0.7462535127781393
$ws.Range("N9").Value = 1.595015869560115
$ws.Range("O9").Value = 1.227727183310407

$ws.Range("B10").Value = 0.1815897941911828
$ws.Range("D10").Value = 0.008392693194828382
$ws.Range("E10").Value = 0.09692274725936301
$ws.Range("F10").Value = 0.461407877260271
$ws.Range("G10").Value = 0.3462889127864912
$ws.Range("H10").Value = 0.36101659559003
$ws.Range("M10").Value = 0.8807565515107427
$ws.Range("N10").Value = 1.704201265293932
$ws.Range("O10").Value = 1.344862573927685

$ws.Range("B11").Value = 0.1899055038935558
$ws.Range("D11").Value = 0.008934375093460289
$ws.Range("E11").Value = 0.09550074769918382
$ws.Range("F11").Value = 0.4806957960630172
$ws.Range("G11").Value = 0.3652480036625008
$ws.Range("H11").Value = 0.367577149910133
$ws.Range("M11").Value = 0.9418041721421844
$ws.Range("N11").Value = 1.754347378327935
$ws.Range("O11").Value = 1.398987853658014

$ws.Range("B12").Value = 0.1930635664858329
$ws.Range("D12").Value = 0.009139513645834541
$ws.Range("E12").Value = 0.09499347649300205
$ws.Range("F12").Value = 0.4880390724165551
$ws.Range("G12").Value = 0.372454003601149
$ws.Range("H12").Value = 0.3700933823506034
$ws.Range("M12").Value = 0.9649006505349291
$ws.Range("N12").Value = 1.773403095966614
$ws.Range("O12").Value = 1.419605371796592

$ws.Range("B13").Value = 0.1923830208949227
$ws.Range("D13").Value = 0.009095332931515543
$ws.Range("E13").Value = 0.0951013316482765
$ws.Range("F13").Value = 0.486455813172185
$ws.Range("G13").Value = 0.3709008749752769
$ws.Range("H13").Value = 0.3695500475313338
$ws.Range("M13").Value = 0.9599273576716314
$ws.Range("N13").Value = 1.769296183641586
$ws.Range("O13").Value = 1.415159607964654

$ws.Range("B14").Value = 0.1901651388805448
$ws.Range("D14").Value = 0.008951251738835708
$ws.Range("E14").Value = 0.09545838628977421
$ws.Range("F14").Value = 0.4812991420065202
$ws.Range("G14").Value = 0.3658403108537982
$ws.Range("H14").Value = 0.3677835224357722
$ws.Range("M14").Value = 0.9437047585989831
$ws.Range("N14").Value = 1.755913783642143
$ws.Range("O14").Value = 1.400681631332304

$ws.Range("B15").Value = 0.1888077981938494
$ws.Range("D15").Value = 0.008862999428195906
$ws.Range("E15").Value = 0.09568116947913552
$ws.Range("F15").Value = 0.4781456629824419
$ws.Range("G15").Value = 0.3627440419991927
$ws.Range("H15").Value = 0.3667056295017233
$ws.Range("M15").Value = 0.9337651924121531
$ws.Range("N15").Value = 1.747725270877936
$ws.Range("O15").Value = 1.391829289236711

$ws.Range("B16").Value = 0.1810476267981187
$ws.Range("D16").Value = 0.008357295935859099
$ws.Range("E16").Value = 0.09702002251664155
$ws.Range("F16").Value = 0.4601528648906879
$ws.Range("G16").Value = 0.3450535993470112
$ws.Range("H16").Value = 0.3605923085063552
$ws.Range("M16").Value = 0.8767640721859635
$ws.Range("N16").Value = 1.700933516733244
$ws.Range("O16").Value = 1.341342317346886

$ws.Range("B17").Value = 0.1763034408433839
$ws.Range("D17").Value = 0.008047105821006539
$ws.Range("E17").Value = 0.09789650079799372
$ws.Range("F17").Value = 0.4491847867966925
$ws.Range("G17").Value = 0.3342481742829193
$ws.Range("H17").Value = 0.3568987427014605
$ws.Range("M17").Value = 0.8417595361178769
$ws.Range("N17").Value = 1.672348967884972
$ws.Range("O17").Value = 1.310585778442629

$ws.Range("B18").Value = 0.1735808223692317
$ws.Range("D18").Value = 0.007868712895270136
$ws.Range("E18").Value = 0.09842075432035102
$ws.Range("F18").Value = 0.4429018783877297
$ws.Range("G18").Value = 0.3280503953542251
$ws.Range("H18").Value = 0.3547951459239442
$ws.Range("M18").Value = 0.8216128781593568
$ws.Range("N18").Value = 1.655952927526073
$ws.Range("O18").Value = 1.292974436068391

$ws.Range("B19").Value = 0.1726600473234328
$ws.Range("D19").Value = 0.007808316021961303
$ws.Range("E19").Value = 0.09860170143116065
$ws.Range("F19").Value = 0.44077899196499
$ws.Range("G19").Value = 0.325954882263801
$ws.Range("H19").Value = 0.3540864829510042
$ws.Range("M19").Value = 0.8147893686155072
$ws.Range("N19").Value = 1.650409310447628
$ws.Range("O19").Value = 1.287025086426922

$ws.Range("B20").Value = 0.1768078366705197
$ws.Range("D20").Value = 0.008080124098782449
$ws.Range("E20").Value = 0.0978011121255058
$ws.Range("F20").Value = 0.4503497025313266
$ws.Range("G20").Value = 0.3353966464481175
$ws.Range("H20").Value = 0.3572897714334857
$ws.Range("M20").Value = 0.8454871789546559
$ws.Range("N20").Value = 1.675387196387021
$ws.Range("O20").Value = 1.313851682775976

$ws.Range("B21").Value = 0.1908163399091478
$ws.Range("D21").Value = 0.008993571573867598
$ws.Range("E21").Value = 0.09535266030610856
$ws.Range("F21").Value = 0.4828127117327909
$ws.Range("G21").Value = 0.3673259969865796
$ws.Range("H21").Value = 0.368301527764757
$ws.Range("M21").Value = 0.9484703053065147
$ws.Range("N21").Value = 1.759842731794549
$ws.Range("O21").Value = 1.404930862457775

$ws.Range("B22").Value = 0.2000245461899652
$ws.Range("D22").Value = 0.009590645302125722
$ws.Range("E22").Value = 0.09393456328208671
$ws.Range("F22").Value = 0.5042586651307062
$ws.Range("G22").Value = 0.3883488677905689
$ws.Range("H22").Value = 0.37568426367298
$ws.Range("M22").Value = 1.015653265042246
$ws.Range("N22").Value = 1.815425783913867
$ws.Range("O22").Value = 1.46516458846645

$ws.Range("B23").Value = 0.1951051951023857
$ws.Range("D23").Value = 0.009271972925660066
$ws.Range("E23").Value = 0.09467462559733875
$ws.Range("F23").Value = 0.4927915000540111
$ws.Range("G23").Value = 0.3771142744199238
$ws.Range("H23").Value = 0.3717269327487145
$ws.Range("M23").Value = 0.9798079948167953
$ws.Range("N23").Value = 1.785725401249294
$ws.Range("O23").Value = 1.432951687095965

$ws.Range("B24").Value = 0.1765797841077728
$ws.Range("D24").Value = 0.008065196722647272
$ws.Range("E24").Value = 0.09784417397323253
$ws.Range("F24").Value = 0.4498229731644869
$ws.Range("G24").Value = 0.3348773773374205
$ws.Range("H24").Value = 0.3571129254090692
$ws.Range("M24").Value = 0.8438019803464698
$ws.Range("N24").Value = 1.67401349600442
$ws.Range("O24").Value = 1.312374946311792

$ws.Range("B25").Value = 0.156831658639689
$ws.Range("D25").Value = 0.006764318950843062
$ws.Range("E25").Value = 0.1020308138487032
$ws.Range("F25").Value = 0.4044538254139809
$ws.Range("G25").Value = 0.2899749218860705
$ws.Range("H25").Value = 0.3421460325277366
$ws.Range("M25").Value = 0.6966050668417978
$ws.Range("N25").Value = 1.555244302609822
$ws.Range("O25").Value = 1.185330592191576
